$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NIDIO")

# Set the File size value for the "nidio_spolis_year_2006_2023" dataset row (row 16)
$ws.Range("E16").Value = "48.0 GB"
